$wb = $excel.ActiveWorkbook

# --- Sheet "Overall Ratings": bump Three Star Rating count from 718 to 722 ---
$wsOverall = $wb.Worksheets.Item("Overall Ratings")
$wsOverall.Range("E2").NumberFormat = "@"
$wsOverall.Range("E2").Value = "722"

# --- Sheet "Employee Reviews": a new review was added at the top, pushing
#     the remaining reviews down by one row and dropping the last (oldest)
#     review off the bottom of the table. ---
$wsReviews = $wb.Worksheets.Item("Employee Reviews")

# Insert a fresh row right below the header, shifting every existing review down.
$wsReviews.Rows.Item(2).Insert()

# The previously last review (row 11 "Field Engineer") is now row 12; remove it
# so the table keeps the same overall size (header + 10 reviews).
$wsReviews.Rows.Item(12).Delete()

# Populate the new review in row 2.
$wsReviews.Range("A2").Value = "Jdss"
$wsReviews.Range("B2").Value = "Bangalore / Bengaluru"
$wsReviews.Range("D2").Value = "Full Time"
$wsReviews.Range("E2").NumberFormat = "@"
$wsReviews.Range("E2").Value = "02 May 2024"
$wsReviews.Range("F2").NumberFormat = "@"
$wsReviews.Range("F2").Value = "2.0"
$wsReviews.Range("G2").Value = "No such thing to like , because company was worst"
$wsReviews.Range("H2").Value = "Without any notice they will fired, hr was not good, not giving fixed salary"
